$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 171: C (nouveaux cas) 7 -> 8, L (nouveaux décès hôpital) 0 -> 1 ---
$ws.Range("C171").Value = 8
$ws.Range("L171").Value = 1

# --- Row 173: C (nouveaux cas) 6 -> 5 ---
$ws.Range("C173").Value = 5

# --- Row 175: C (nouveaux cas) 9 -> 10 ---
$ws.Range("C175").Value = 10

# --- Row 176: C (nouveaux cas) 1 -> 3 ---
$ws.Range("C176").Value = 3

# --- Row 177: C (nouveaux cas) 0 -> 5, D (nouvelles admissions) 0 -> 2 ---
$ws.Range("C177").Value = 5
$ws.Range("D177").Value = 2

# --- Row 178: fill in the day's figures (previously blank) ---
$ws.Range("C178").Value = 0
$ws.Range("D178").Value = 0
$ws.Range("E178").Value = 1
$ws.Range("F178").Value = 1
$ws.Range("G178").Value = 7
$ws.Range("I178").Value = 0
$ws.Range("L178").Value = 0
$ws.Range("M178").Value = 0

# --- Update frozen-pane scroll position / active selection to match the new data entry point ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 164
$excel.ActiveWindow.ScrollColumn = 8
$ws.Range("L178").Select()
